$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "dbsize" table (rows 1-7) originally has 6 experiment column-groups:
#   B=fwd-100  H=bwd-100  N=fwd-1000  T=bwd-1000  Z=fwd-10000  AF=bwd-10000
# A new pair of column-groups (fwd-5000 / bwd-5000) is inserted right after
# bwd-1000 (cols T:Y), so the existing fwd-10000 / bwd-10000 data (cols
# Z:AK) is shifted 12 columns right to AL:AW, and the new fwd-5000 /
# bwd-5000 data takes over the now-vacated Z:AK columns.
# ---------------------------------------------------------------------------

# --- Step 1: Row 2 headers -- move the old group labels out of the way ----
$ws.Range("AL2").Value = $ws.Range("Z2").Value()    # "fwd-10000"
$ws.Range("AR2").Value = $ws.Range("AF2").Value()   # "bwd-10000"

# --- Step 2: Row 3 index row (1,2,3,4,5,avg) is reused verbatim -----------
for ($c = 26; $c -le 37; $c++) {
    $destCol = $c + 12
    $ws.Cells.Item(3, $destCol).Value = $ws.Cells.Item(3, $c).Value()
}

# --- Step 3: Rows 4-7 -- move the raw fwd-10000 / bwd-10000 data ----------
for ($r = 4; $r -le 7; $r++) {
    for ($c = 26; $c -le 37; $c++) {
        $destCol = $c + 12
        $ws.Cells.Item($r, $destCol).Value = $ws.Cells.Item($r, $c).Value()
    }
}

# Re-point the moved average formulas (now in AQ/AW) at their new ranges
$ws.Range("AQ4").Formula = "=AVERAGE(AL4:AP4)"
$ws.Range("AQ5").Formula = "=AVERAGE(AL5:AP5)"
$ws.Range("AQ6").Formula = "=AVERAGE(AL6:AP6)"
$ws.Range("AQ7").Formula = "=AVERAGE(AL7:AP7)"
$ws.Range("AW4").Formula = "=AVERAGE(AR4:AV4)"
$ws.Range("AW5").Formula = "=AVERAGE(AR5:AV5)"
$ws.Range("AW6").Formula = "=AVERAGE(AR6:AV6)"
$ws.Range("AW7").Formula = "=AVERAGE(AR7:AV7)"

# --- Step 4: new headers for the inserted fwd-5000 / bwd-5000 groups ------
$ws.Range("Z2").Value = "fwd-5000"
$ws.Range("AF2").Value = "bwd-5000"

# --- Step 5: new fwd-5000 raw data (cols Z:AD) + average (AE) -------------
$ws.Range("Z4").Value = 0.007326
$ws.Range("AA4").Value = 0.006753
$ws.Range("AB4").Value = 0.006863
$ws.Range("AC4").Value = 0.007153
$ws.Range("AD4").Value = 0.010432
$ws.Range("AE4").Formula = "=AVERAGE(Z4:AD4)"

$ws.Range("Z5").Value = 0.004656
$ws.Range("AA5").Value = 0.004424
$ws.Range("AB5").Value = 0.004185
$ws.Range("AC5").Value = 0.004112
$ws.Range("AD5").Value = 0.004493
$ws.Range("AE5").Formula = "=AVERAGE(Z5:AD5)"

$ws.Range("Z6").Value = 0.028779
$ws.Range("AA6").Value = 0.025654
$ws.Range("AB6").Value = 0.029053
$ws.Range("AC6").Value = 0.026107
$ws.Range("AD6").Value = 0.029322
$ws.Range("AE6").Formula = "=AVERAGE(Z6:AD6)"

$ws.Range("Z7").Value = 0.024795
$ws.Range("AA7").Value = 0.023858
$ws.Range("AB7").Value = 0.022079
$ws.Range("AC7").Value = 0.021487
$ws.Range("AD7").Value = 0.023858
$ws.Range("AE7").Formula = "=AVERAGE(Z7:AD7)"

# --- Step 6: new bwd-5000 raw data (cols AF:AJ) + average (AK) ------------
$ws.Range("AF4").Value = 0.006818
$ws.Range("AG4").Value = 0.006412
$ws.Range("AH4").Value = 0.006647
$ws.Range("AI4").Value = 0.006351
$ws.Range("AJ4").Value = 0.009153
$ws.Range("AK4").Formula = "=AVERAGE(AF4:AJ4)"

$ws.Range("AF5").Value = 0.014656
$ws.Range("AG5").Value = 0.01024
$ws.Range("AH5").Value = 0.014091
$ws.Range("AI5").Value = 0.01031
$ws.Range("AJ5").Value = 0.013845
$ws.Range("AK5").Formula = "=AVERAGE(AF5:AJ5)"

$ws.Range("AF6").Value = 0.037708
$ws.Range("AG6").Value = 0.035338
$ws.Range("AH6").Value = 0.039656
$ws.Range("AI6").Value = 0.038645
$ws.Range("AJ6").Value = 0.039337
$ws.Range("AK6").Formula = "=AVERAGE(AF6:AJ6)"

$ws.Range("AF7").Value = 0.021187
$ws.Range("AG7").Value = 0.021586
$ws.Range("AH7").Value = 0.02188
$ws.Range("AI7").Value = 0.0257
$ws.Range("AJ7").Value = 0.024653
$ws.Range("AK7").Formula = "=AVERAGE(AF7:AJ7)"

# --- Step 7: sheet view cosmetics (match the saved selection) -------------
[void]$ws.Range("AK7").Select()

Write-Output "edit applied"
